{"js": "// 1) Remove the trailing \", \" that follows \"Laura Moreno Ramos\" in the\n//    \"D/D\u00aa ...\" paragraph (the whole run containing just \", \" is deleted).\nconst nameResults = context.document.body.search(\"Ramos, \", { matchCase: true });\nnameResults.load(\"text\");\nawait context.sync();\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"Ramos\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Tidy up \"Desarrollo de Aplicaciones Web \" -> \"Desarrollo de Aplicaciones Web\"\n//    (drop the trailing space that duplicated the following run's leading space).\nconst degreeResults = context.document.body.search(\"Desarrollo de Aplicaciones Web \", { matchCase: true });\ndegreeResults.load(\"text\");\nawait context.sync();\nif (degreeResults.items.length > 0) {\n  degreeResults.items[0].insertText(\"Desarrollo de Aplicaciones Web\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Collapse the run \" de        \" (multiple spaces) down to \" de \".\nconst deResults = context.document.body.search(\" de        \", { matchCase: true });\ndeResults.load(\"text\");\nawait context.sync();\nif (deResults.items.length > 0) {\n  deResults.items[0].insertText(\" de \", \"Replace\");\n  await context.sync();\n}\n\n// 4) Update the signature date: \"29\" -> \"24\" and \"abril\" -> \"mayo\".\nconst dayResults = context.document.body.search(\"29\", { matchCase: true });\ndayResults.load(\"text\");\nawait context.sync();\nif (dayResults.items.length > 0) {\n  dayResults.items[0].insertText(\"24\", \"Replace\");\n  await context.sync();\n}\n\nconst monthResults = context.document.body.search(\"abril\", { matchCase: true });\nmonthResults.load(\"text\");\nawait context.sync();\nif (monthResults.items.length > 0) {\n  monthResults.items[0].insertText(\"mayo\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the trailing \", \" that follows \"Laura Moreno Ramos\" in the\n#    \"D/D\u00aa ...\" paragraph (the whole run containing just \", \" is deleted).\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"Ramos, \"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Ramos\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 1) | Out-Null\n\n# 2) Tidy up \"Desarrollo de Aplicaciones Web \" -> \"Desarrollo de Aplicaciones Web\"\n#    (drop the trailing space that duplicated the following run's leading space).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Desarrollo de Aplicaciones Web \"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Desarrollo de Aplicaciones Web\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 1) | Out-Null\n\n# 3) Collapse the run \" de        \" (multiple spaces) down to \" de \".\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \" de        \"\n$find3.Replacement.ClearFormatting()\n$find3.Replacement.Text = \" de \"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 1) | Out-Null\n\n# 4) Update the signature date: \"29\" -> \"24\" and \"abril\" -> \"mayo\".\n$find4 = $d.Content.Find\n$find4.ClearFormatting()\n$find4.Text = \"29\"\n$find4.Replacement.ClearFormatting()\n$find4.Replacement.Text = \"24\"\n$find4.Execute($find4.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find4.Replacement.Text, 1) | Out-Null\n\n$find5 = $d.Content.Find\n$find5.ClearFormatting()\n$find5.Text = \"abril\"\n$find5.Replacement.ClearFormatting()\n$find5.Replacement.Text = \"mayo\"\n$find5.Execute($find5.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find5.Replacement.Text, 1) | Out-Null\n"}
